# Insert a new data row at row 61 (pushes existing rows 61-157 down to 62-158),
# then populate the new row with its values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(61).Insert()

$ws.Range("A61").Value = 11
$ws.Range("B61").Value = "Vega Monumental Concepción"
$ws.Range("C61").Value = "Bíobío"
$ws.Range("D61").Value = 44848
$ws.Range("E61").Value = 8
$ws.Range("F61").Value = 100112032
$ws.Range("G61").Value = "Zapallo italiano"
$ws.Range("H61").Value = "Sin especificar"
$ws.Range("I61").Value = "Primera"
$ws.Range("J61").Value = 310
$ws.Range("K61").Value = 17000
$ws.Range("L61").Value = 19000
$ws.Range("M61").Value = 17968
$ws.Range("N61").Value = "`$/caja 50 unidades"
$ws.Range("O61").Value = "Región de O'Higgins"
$ws.Range("P61").Value = 359
$ws.Range("Q61").Value = 50
$ws.Range("R61").Value = "Hortaliza"
